# Generate Report for Handback
# Update the timestamp strings recorded in the handback status report.
# The cells already carry a text-style number format (yyyy-mm-dd HH:mm:ss)
# but store their content as literal shared-string text, so we assign the
# new value as text to keep them stored the same way.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview!G2 - Latest HO Xliff Generate Date
$overview.Range("G2").Value = "2016-08-31 17:22:16"

# zh-cn!H2 - Correspond Handoff Datetime
$zhcn.Range("H2").Value = "2016-08-31 17:22:05"

# zh-cn!K2 - Correspond Handback DateTime
$zhcn.Range("K2").Value = "2016-08-31 17:22:31"

# de-de!H2 - Correspond Handoff Datetime (shares the same timestamp text as Overview!G2)
$dede.Range("H2").Value = "2016-08-31 17:22:16"

# de-de!K2 - Correspond Handback DateTime
$dede.Range("K2").Value = "2016-08-31 17:22:38"
